# NFP / Initial Jobless Claims data preprocessing refactor
# Updates the Actual/Revised figures, summary statistics, a few economist
# ranks, and appends a new economist estimate row (Marc Giannoni / Barclays).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Actual / Revised (previously placeholder "nan" text values) ---
$ws.Range("C7").Value = 221
$ws.Range("C9").Value = 228

# --- Summary block ---
$ws.Range("C12").Value = "233k"
$ws.Range("C13").Value = "233.01k"
$ws.Range("C16").Value = 43
$ws.Range("C18").Value = "4.37k"

# --- Rank corrections among existing economist rows ---
$ws.Range("F28").Value = "4th"
$ws.Range("F34").Value = "5th"
$ws.Range("F48").Value = "6th"
$ws.Range("F54").Value = ""

# --- Append a new economist estimate row (row 65), push the trailing
#     blank spacer row down to row 66 ---
$ws.Range("B65").Value = "Marc Giannoni"
$ws.Range("C65").Value = "Barclays Capital Inc."
$ws.Range("D65").Value = 230
$ws.Range("E65").Value = "7/16/2025"

$ws.Range("A66:F66").Value = ""
